$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Characters(21,2).Text = "16"
$ws.Range("C9").Characters(27,9).Text = "4/17/2023"
$ws.Range("C9").Characters(47,9).Text = "4/23/2023"

# --- Fix cell formatting for cells whose type flips between number <-> text. ---
# Copying a same-column donor cell brings over both the xf style AND the
# shared-string/number typing in one shot; for cells becoming the "0" / "***.*"
# sentinel text, the donor value (row 23, which is already all sentinels) is
# exactly the target, so nothing more is needed. For cells becoming plain
# numbers, the donor (row 16, all plain numbers) just seeds the numeric style;
# the real value is written below.
$ws.Range("C23").Copy($ws.Range("C14"))
$ws.Range("G23").Copy($ws.Range("G14"))
$ws.Range("H23").Copy($ws.Range("H14"))
$ws.Range("D16").Copy($ws.Range("D15"))
$ws.Range("E16").Copy($ws.Range("E15"))
$ws.Range("D23").Copy($ws.Range("D22"))
$ws.Range("E23").Copy($ws.Range("E22"))
$ws.Range("D16").Copy($ws.Range("D26"))
$ws.Range("E16").Copy($ws.Range("E26"))
$ws.Range("F16").Copy($ws.Range("F26"))
$ws.Range("C16").Copy($ws.Range("C27"))
$ws.Range("G23").Copy($ws.Range("G28"))
$ws.Range("H23").Copy($ws.Range("H28"))
$ws.Range("G23").Copy($ws.Range("G29"))
$ws.Range("H23").Copy($ws.Range("H29"))
$ws.Range("D16").Copy($ws.Range("D30"))
$ws.Range("E16").Copy($ws.Range("E30"))
$ws.Range("G16").Copy($ws.Range("G30"))
$ws.Range("H16").Copy($ws.Range("H30"))

# --- Updated crime-stat grid values ---
$values = [ordered]@{
  "D15" = 1
  "E15" = -100
  "G15" = 2
  "J15" = 8
  "K15" = -12.5
  "C16" = 1
  "E16" = -50
  "F16" = 12
  "G16" = 9
  "H16" = 33.333333333333
  "I16" = 38
  "J16" = 32
  "K16" = 18.75
  "L16" = 46.153846153846
  "M16" = -28.301886792452
  "N16" = -86.713286713286
  "C17" = 3
  "D17" = 7
  "E17" = -57.142857142857
  "F17" = 11
  "G17" = 16
  "H17" = -31.25
  "I17" = 61
  "J17" = 40
  "K17" = 52.5
  "L17" = 5.172413793103
  "M17" = 90.625
  "N17" = -24.691358024691
  "C18" = 1
  "D18" = 2
  "E18" = -50
  "F18" = 10
  "G18" = 10
  "H18" = 0
  "I18" = 81
  "J18" = 53
  "K18" = 52.830188679245
  "L18" = 84.090909090909
  "M18" = -11.956521739130
  "N18" = -82.543103448275
  "C19" = 19
  "D19" = 19
  "E19" = 0
  "F19" = 50
  "G19" = 48
  "H19" = 4.166666666666
  "I19" = 193
  "J19" = 223
  "K19" = -13.452914798206
  "L19" = 19.875776397515
  "M19" = 55.645161290322
  "N19" = -13.063063063063
  "C20" = 5
  "D20" = 1
  "E20" = 400
  "F20" = 16
  "G20" = 4
  "H20" = 300
  "I20" = 46
  "J20" = 27
  "K20" = 70.370370370370
  "L20" = 155.555555555556
  "M20" = -23.333333333333
  "N20" = -91.958041958042
  "C21" = 29
  "D21" = 32
  "E21" = -9.375
  "F21" = 100
  "G21" = 89
  "H21" = 12.359550561797
  "I21" = 427
  "J21" = 386
  "K21" = 10.621761658031
  "L21" = 37.299035369774
  "M21" = 15.718157181571
  "N21" = -73.883792048929
  "M22" = -77.777777777777
  "C24" = 40
  "D24" = 26
  "E24" = 53.846153846153
  "F24" = 145
  "G24" = 121
  "H24" = 19.834710743801
  "I24" = 533
  "J24" = 489
  "K24" = 8.997955010224
  "L24" = 61.027190332326
  "M24" = 89.007092198581
  "C25" = 14
  "D25" = 5
  "E25" = 180
  "F25" = 44
  "G25" = 25
  "H25" = 76
  "I25" = 167
  "J25" = 129
  "K25" = 29.457364341085
  "L25" = 31.496062992126
  "M25" = 36.885245901639
  "D26" = 1
  "E26" = -100
  "F26" = 1
  "G26" = 2
  "H26" = -50
  "I26" = 8
  "J26" = 11
  "K26" = -27.272727272727
  "L26" = 100
  "C27" = 2
  "E27" = 0
  "F27" = 4
  "H27" = -42.857142857142
  "I27" = 10
  "J27" = 15
  "K27" = -33.333333333333
  "L27" = -16.666666666666
  "D30" = 1
  "E30" = -100
  "G30" = 1
  "H30" = 0
  "J30" = 5
  "K30" = -80
  "L30" = -75
}
foreach ($ref in $values.Keys) {
  $ws.Range($ref).Value = $values[$ref]
}
